# Auto-generated by analysis script
$wb = $excel.ActiveWorkbook
$wsVars = $wb.Worksheets.Item(1)
$wsCats = $wb.Worksheets.Item(2)

# --- Variables sheet: insert 3 rows at row 5 to make room for the 4 new parental-education rows ---
$wsVars.Range("A5:A7").EntireRow.Insert()

# Fill rows 2..57 of the Variables sheet with the final data
$wsVars.Cells.Item(2, 1).Value = 1
$wsVars.Cells.Item(2, 2).Value = "ID"
$wsVars.Cells.Item(2, 3).Value = "ID"
$wsVars.Cells.Item(2, 4).Value = "integer"
$wsVars.Cells.Item(3, 1).Value = 2
$wsVars.Cells.Item(3, 2).Value = "sex"
$wsVars.Cells.Item(3, 3).Value = "Geschlecht"
$wsVars.Cells.Item(3, 4).Value = "integer"
$wsVars.Cells.Item(4, 1).Value = 3
$wsVars.Cells.Item(4, 2).Value = "AGE_BASE"
$wsVars.Cells.Item(4, 3).Value = "Age [years] at dietary assessment (1. record day)"
$wsVars.Cells.Item(4, 4).Value = "decimal"
$wsVars.Cells.Item(5, 1).Value = 4
$wsVars.Cells.Item(5, 2).Value = "m_schulab"
$wsVars.Cells.Item(5, 3).Value = "Schuldbildung der Mutter"
$wsVars.Cells.Item(5, 4).Value = "integer"
$wsVars.Cells.Item(6, 1).Value = 5
$wsVars.Cells.Item(6, 2).Value = "v_schulab"
$wsVars.Cells.Item(6, 3).Value = "Schuldbildung des Vaters"
$wsVars.Cells.Item(6, 4).Value = "integer"
$wsVars.Cells.Item(7, 1).Value = 6
$wsVars.Cells.Item(7, 2).Value = "b_berufsab"
$wsVars.Cells.Item(7, 3).Value = "Berufsabschluss der Mutter"
$wsVars.Cells.Item(7, 4).Value = "integer"
$wsVars.Cells.Item(8, 1).Value = 7
$wsVars.Cells.Item(8, 2).Value = "v_berufsab"
$wsVars.Cells.Item(8, 3).Value = "Berufsabschluss des Vaters"
$wsVars.Cells.Item(8, 4).Value = "integer"
$wsVars.Cells.Item(9, 1).Value = 8
$wsVars.Cells.Item(9, 2).Value = "EMPLOY_P"
$wsVars.Cells.Item(9, 3).Value = "parental employment status"
$wsVars.Cells.Item(9, 4).Value = "integer"
$wsVars.Cells.Item(10, 1).Value = 9
$wsVars.Cells.Item(10, 2).Value = "TOT_PA_QX"
$wsVars.Cells.Item(10, 3).Value = "Physical activity from questionnaire data [MET-hr/week]"
$wsVars.Cells.Item(10, 4).Value = "decimal"
$wsVars.Cells.Item(11, 1).Value = 10
$wsVars.Cells.Item(11, 2).Value = "SMOKE_ST_HH"
$wsVars.Cells.Item(11, 3).Value = "Smoker in Household [yes/no]"
$wsVars.Cells.Item(11, 4).Value = "integer"
$wsVars.Cells.Item(12, 1).Value = 11
$wsVars.Cells.Item(12, 2).Value = "TG"
$wsVars.Cells.Item(12, 3).Value = "Triglycerides measured from blood samples [mg/dl]"
$wsVars.Cells.Item(12, 4).Value = "decimal"
$wsVars.Cells.Item(13, 1).Value = 12
$wsVars.Cells.Item(13, 2).Value = "CHOL"
$wsVars.Cells.Item(13, 3).Value = "Cholesterol measured from blood samples [mg/dl]"
$wsVars.Cells.Item(13, 4).Value = "decimal"
$wsVars.Cells.Item(14, 1).Value = 13
$wsVars.Cells.Item(14, 2).Value = "LDL"
$wsVars.Cells.Item(14, 3).Value = "LDL measured from blood samples [mg/dl]"
$wsVars.Cells.Item(14, 4).Value = "decimal"
$wsVars.Cells.Item(15, 1).Value = 14
$wsVars.Cells.Item(15, 2).Value = "HDL"
$wsVars.Cells.Item(15, 3).Value = "HDL measured from blood samples [mg/dl]"
$wsVars.Cells.Item(15, 4).Value = "decimal"
$wsVars.Cells.Item(16, 1).Value = 15
$wsVars.Cells.Item(16, 2).Value = "BMI"
$wsVars.Cells.Item(16, 3).Value = "BMI at dietary assessment"
$wsVars.Cells.Item(16, 4).Value = "decimal"
$wsVars.Cells.Item(17, 1).Value = 16
$wsVars.Cells.Item(17, 2).Value = "BMI_FUP"
$wsVars.Cells.Item(17, 3).Value = "BMI at follow-up"
$wsVars.Cells.Item(17, 4).Value = "decimal"
$wsVars.Cells.Item(18, 1).Value = 17
$wsVars.Cells.Item(18, 2).Value = "BMI_SDS"
$wsVars.Cells.Item(18, 3).Value = "BMI-SDS at dietary assessment"
$wsVars.Cells.Item(18, 4).Value = "decimal"
$wsVars.Cells.Item(19, 1).Value = 18
$wsVars.Cells.Item(19, 2).Value = "WAIST_FUP"
$wsVars.Cells.Item(19, 3).Value = "Waist circumference at follow-up [cm]"
$wsVars.Cells.Item(19, 4).Value = "decimal"
$wsVars.Cells.Item(20, 1).Value = 19
$wsVars.Cells.Item(20, 2).Value = "WAIST"
$wsVars.Cells.Item(20, 3).Value = "Waist circumference at dietary assessment [cm]"
$wsVars.Cells.Item(20, 4).Value = "decimal"
$wsVars.Cells.Item(21, 1).Value = 20
$wsVars.Cells.Item(21, 2).Value = "HIP"
$wsVars.Cells.Item(21, 3).Value = "Hip circumference at dietary assessment [cm]"
$wsVars.Cells.Item(21, 4).Value = "decimal"
$wsVars.Cells.Item(22, 1).Value = 21
$wsVars.Cells.Item(22, 2).Value = "HIP_FUP"
$wsVars.Cells.Item(22, 3).Value = "Hip circumference at follow-up [cm]"
$wsVars.Cells.Item(22, 4).Value = "decimal"
$wsVars.Cells.Item(23, 1).Value = 22
$wsVars.Cells.Item(23, 2).Value = "FMI_FUP"
$wsVars.Cells.Item(23, 3).Value = "FMI at follow-up"
$wsVars.Cells.Item(23, 4).Value = "decimal"
$wsVars.Cells.Item(24, 1).Value = 23
$wsVars.Cells.Item(24, 2).Value = "FMI"
$wsVars.Cells.Item(24, 3).Value = "FMI at dietary assessment"
$wsVars.Cells.Item(24, 4).Value = "decimal"
$wsVars.Cells.Item(25, 1).Value = 24
$wsVars.Cells.Item(25, 2).Value = "BODY_FAT_FUP"
$wsVars.Cells.Item(25, 3).Value = "Body fat % according to Durnin and Womersley at FUP"
$wsVars.Cells.Item(25, 4).Value = "decimal"
$wsVars.Cells.Item(26, 1).Value = 25
$wsVars.Cells.Item(26, 2).Value = "BODY_FAT"
$wsVars.Cells.Item(26, 3).Value = "Body fat % according to Slaughter at dietary assessment"
$wsVars.Cells.Item(26, 4).Value = "decimal"
$wsVars.Cells.Item(27, 1).Value = 26
$wsVars.Cells.Item(27, 2).Value = "AGE_ANTH_FUP"
$wsVars.Cells.Item(27, 3).Value = "age at anthro measurement"
$wsVars.Cells.Item(27, 4).Value = "decimal"
$wsVars.Cells.Item(28, 1).Value = 27
$wsVars.Cells.Item(28, 2).Value = "ENERGY"
$wsVars.Cells.Item(28, 3).Value = "Energy intake [kcal/d]"
$wsVars.Cells.Item(28, 4).Value = "decimal"
$wsVars.Cells.Item(29, 1).Value = 28
$wsVars.Cells.Item(29, 2).Value = "CARB"
$wsVars.Cells.Item(29, 3).Value = "Carbohydrate intake [g/d]"
$wsVars.Cells.Item(29, 4).Value = "decimal"
$wsVars.Cells.Item(30, 1).Value = 29
$wsVars.Cells.Item(30, 2).Value = "PROT"
$wsVars.Cells.Item(30, 3).Value = "Protein intake [g/d]"
$wsVars.Cells.Item(30, 4).Value = "decimal"
$wsVars.Cells.Item(31, 1).Value = 30
$wsVars.Cells.Item(31, 2).Value = "FAT"
$wsVars.Cells.Item(31, 3).Value = "Fat intake [g/d]"
$wsVars.Cells.Item(31, 4).Value = "decimal"
$wsVars.Cells.Item(32, 1).Value = 31
$wsVars.Cells.Item(32, 2).Value = "ALC"
$wsVars.Cells.Item(32, 3).Value = "Alcohol intake [g/d]"
$wsVars.Cells.Item(32, 4).Value = "decimal"
$wsVars.Cells.Item(33, 1).Value = 32
$wsVars.Cells.Item(33, 2).Value = "FIBER"
$wsVars.Cells.Item(33, 3).Value = "Dietary fiber intake [g/d]"
$wsVars.Cells.Item(33, 4).Value = "decimal"
$wsVars.Cells.Item(34, 1).Value = 33
$wsVars.Cells.Item(34, 2).Value = "SFA"
$wsVars.Cells.Item(34, 3).Value = "Saturated fat intake [g/d]"
$wsVars.Cells.Item(34, 4).Value = "decimal"
$wsVars.Cells.Item(35, 1).Value = 34
$wsVars.Cells.Item(35, 2).Value = "MUFA"
$wsVars.Cells.Item(35, 3).Value = "Monounsaturated fat intake [g/d]"
$wsVars.Cells.Item(35, 4).Value = "decimal"
$wsVars.Cells.Item(36, 1).Value = 35
$wsVars.Cells.Item(36, 2).Value = "PUFA"
$wsVars.Cells.Item(36, 3).Value = "Polyunsaturated fat intake [g/d]"
$wsVars.Cells.Item(36, 4).Value = "decimal"
$wsVars.Cells.Item(37, 1).Value = 36
$wsVars.Cells.Item(37, 2).Value = "TOT_SUGARS"
$wsVars.Cells.Item(37, 3).Value = "Total sugar intake [g/d]"
$wsVars.Cells.Item(37, 4).Value = "decimal"
$wsVars.Cells.Item(38, 1).Value = 37
$wsVars.Cells.Item(38, 2).Value = "ADD_SUGARS"
$wsVars.Cells.Item(38, 3).Value = "Added sugar intake [g/d]"
$wsVars.Cells.Item(38, 4).Value = "decimal"
$wsVars.Cells.Item(39, 1).Value = 38
$wsVars.Cells.Item(39, 2).Value = "FREE_SUGARS"
$wsVars.Cells.Item(39, 3).Value = "Free sugar intake [g/d]"
$wsVars.Cells.Item(39, 4).Value = "decimal"
$wsVars.Cells.Item(40, 1).Value = 39
$wsVars.Cells.Item(40, 2).Value = "GLUC"
$wsVars.Cells.Item(40, 3).Value = "Glucose intake [g/d]"
$wsVars.Cells.Item(40, 4).Value = "decimal"
$wsVars.Cells.Item(41, 1).Value = 40
$wsVars.Cells.Item(41, 2).Value = "FRUC"
$wsVars.Cells.Item(41, 3).Value = "Fructose intake [g/d]"
$wsVars.Cells.Item(41, 4).Value = "decimal"
$wsVars.Cells.Item(42, 1).Value = 41
$wsVars.Cells.Item(42, 2).Value = "GI"
$wsVars.Cells.Item(42, 3).Value = "Daily glycaemic index"
$wsVars.Cells.Item(42, 4).Value = "decimal"
$wsVars.Cells.Item(43, 1).Value = 42
$wsVars.Cells.Item(43, 2).Value = "GL"
$wsVars.Cells.Item(43, 3).Value = "Daily glycaemic load"
$wsVars.Cells.Item(43, 4).Value = "decimal"
$wsVars.Cells.Item(44, 1).Value = 43
$wsVars.Cells.Item(44, 2).Value = "SODIUM"
$wsVars.Cells.Item(44, 3).Value = "Sodium intake [mg/d]"
$wsVars.Cells.Item(44, 4).Value = "decimal"
$wsVars.Cells.Item(45, 1).Value = 44
$wsVars.Cells.Item(45, 2).Value = "SOD_POT_RATIO"
$wsVars.Cells.Item(45, 3).Value = "Sodium to potassium intake ratio [g/d]"
$wsVars.Cells.Item(45, 4).Value = "decimal"
$wsVars.Cells.Item(46, 1).Value = 45
$wsVars.Cells.Item(46, 2).Value = "SUGAR_CONFECT_11"
$wsVars.Cells.Item(46, 3).Value = "Intake of sugar and similar, confectionery and water-based sweet desserts [g/d]"
$wsVars.Cells.Item(46, 4).Value = "decimal"
$wsVars.Cells.Item(47, 1).Value = 46
$wsVars.Cells.Item(47, 2).Value = "CAKES_12"
$wsVars.Cells.Item(47, 3).Value = "Intake of cakes and fine bakery products [g/d]"
$wsVars.Cells.Item(47, 4).Value = "decimal"
$wsVars.Cells.Item(48, 1).Value = 47
$wsVars.Cells.Item(48, 2).Value = "FRUITVEG_JUICE_1301"
$wsVars.Cells.Item(48, 3).Value = "Intake of fruit and vegetable juices [g/d]"
$wsVars.Cells.Item(48, 4).Value = "decimal"
$wsVars.Cells.Item(49, 1).Value = 48
$wsVars.Cells.Item(49, 2).Value = "SOFTDRINKS_1302"
$wsVars.Cells.Item(49, 3).Value = "Intake of soft drinks [g/d]"
$wsVars.Cells.Item(49, 4).Value = "decimal"
$wsVars.Cells.Item(50, 1).Value = 49
$wsVars.Cells.Item(50, 2).Value = "ART_SWEETENER_170201"
$wsVars.Cells.Item(50, 3).Value = "Intake of artificial sweeteners (e.g., aspartam, saccharine) [g/d]"
$wsVars.Cells.Item(50, 4).Value = "decimal"
$wsVars.Cells.Item(51, 1).Value = 50
$wsVars.Cells.Item(51, 2).Value = "VEGETABLES_02"
$wsVars.Cells.Item(51, 3).Value = "Vegetable intake [g/d]"
$wsVars.Cells.Item(51, 4).Value = "decimal"
$wsVars.Cells.Item(52, 1).Value = 51
$wsVars.Cells.Item(52, 2).Value = "LEGUMES_TOT_03"
$wsVars.Cells.Item(52, 3).Value = "Total legumes intake [g/d]"
$wsVars.Cells.Item(52, 4).Value = "decimal"
$wsVars.Cells.Item(53, 1).Value = 52
$wsVars.Cells.Item(53, 2).Value = "FRUITS_TOT_04"
$wsVars.Cells.Item(53, 3).Value = "Total fruit intake [g/d]"
$wsVars.Cells.Item(53, 4).Value = "decimal"
$wsVars.Cells.Item(54, 1).Value = 53
$wsVars.Cells.Item(54, 2).Value = "RED_MEAT_0701"
$wsVars.Cells.Item(54, 3).Value = "Intake of red meat (mammals meat) [g/d]"
$wsVars.Cells.Item(54, 4).Value = "decimal"
$wsVars.Cells.Item(55, 1).Value = 54
$wsVars.Cells.Item(55, 2).Value = "PROCMEAT_0704"
$wsVars.Cells.Item(55, 3).Value = "Intake of processed or preserved meat [g/d]"
$wsVars.Cells.Item(55, 4).Value = "decimal"
$wsVars.Cells.Item(56, 1).Value = 55
$wsVars.Cells.Item(56, 2).Value = "COFFEE_130301"
$wsVars.Cells.Item(56, 3).Value = "Coffee intake [g/d]"
$wsVars.Cells.Item(56, 4).Value = "decimal"
$wsVars.Cells.Item(57, 1).Value = 56
$wsVars.Cells.Item(57, 2).Value = "TEA_130302"
$wsVars.Cells.Item(57, 3).Value = "Tea intake [g/d]"
$wsVars.Cells.Item(57, 4).Value = "decimal"

# --- Categories sheet: remove old EDU_LEVEL_P (10 rows) + EMPLOY_P (7 rows) category rows (4..20), then insert 26 blank rows back ---
$wsCats.Range("A4:A20").EntireRow.Delete()
$wsCats.Range("A4:A29").EntireRow.Insert()

# Fill rows 4..29 with the new category rows (m_schulab, v_schulab, m_berufsab, v_berufsab)
$wsCats.Cells.Item(4, 1).Value = "m_schulab"
$wsCats.Cells.Item(4, 2).Value = 1
$wsCats.Cells.Item(4, 3).Value = "Volks-/Hauptschulabschluss"
$wsCats.Cells.Item(5, 1).Value = "m_schulab"
$wsCats.Cells.Item(5, 2).Value = 2
$wsCats.Cells.Item(5, 3).Value = "Mittlere Reife, Realschulabschluss (Fachschulreife)"
$wsCats.Cells.Item(6, 1).Value = "m_schulab"
$wsCats.Cells.Item(6, 2).Value = 3
$wsCats.Cells.Item(6, 3).Value = "Fachhochschulreife, Abschluss einer Fachoberschule etc.       "
$wsCats.Cells.Item(7, 1).Value = "m_schulab"
$wsCats.Cells.Item(7, 2).Value = 4
$wsCats.Cells.Item(7, 3).Value = "Abitur (Hochschulreife)"
$wsCats.Cells.Item(8, 1).Value = "m_schulab"
$wsCats.Cells.Item(8, 2).Value = 5
$wsCats.Cells.Item(8, 3).Value = "keinen dieser Abschlüsse"
$wsCats.Cells.Item(9, 1).Value = "v_schulab"
$wsCats.Cells.Item(9, 2).Value = 1
$wsCats.Cells.Item(9, 3).Value = "Volks-/Hauptschulabschluss"
$wsCats.Cells.Item(10, 1).Value = "v_schulab"
$wsCats.Cells.Item(10, 2).Value = 2
$wsCats.Cells.Item(10, 3).Value = "Mittlere Reife, Realschulabschluss (Fachschulreife)"
$wsCats.Cells.Item(11, 1).Value = "v_schulab"
$wsCats.Cells.Item(11, 2).Value = 3
$wsCats.Cells.Item(11, 3).Value = "Fachhochschulreife, Abschluss einer Fachoberschule etc.       "
$wsCats.Cells.Item(12, 1).Value = "v_schulab"
$wsCats.Cells.Item(12, 2).Value = 4
$wsCats.Cells.Item(12, 3).Value = "Abitur (Hochschulreife)"
$wsCats.Cells.Item(13, 1).Value = "v_schulab"
$wsCats.Cells.Item(13, 2).Value = 5
$wsCats.Cells.Item(13, 3).Value = "keinen dieser Abschlüsse"
$wsCats.Cells.Item(14, 1).Value = "m_berufsab"
$wsCats.Cells.Item(14, 2).Value = 1
$wsCats.Cells.Item(14, 3).Value = "Lehre (beruflich-betriebliche Ausbildung)"
$wsCats.Cells.Item(15, 1).Value = "m_berufsab"
$wsCats.Cells.Item(15, 2).Value = 2
$wsCats.Cells.Item(15, 3).Value = "Berufsschule, Handelsschule (berufl.-schulische Ausbildung)"
$wsCats.Cells.Item(16, 1).Value = "m_berufsab"
$wsCats.Cells.Item(16, 2).Value = 3
$wsCats.Cells.Item(16, 3).Value = "Fachschule (z.B. Meister-Technikerschule, Berufs/Fachakademie)"
$wsCats.Cells.Item(17, 1).Value = "m_berufsab"
$wsCats.Cells.Item(17, 2).Value = 4
$wsCats.Cells.Item(17, 3).Value = "Fachhochschule, Ingenieurschule"
$wsCats.Cells.Item(18, 1).Value = "m_berufsab"
$wsCats.Cells.Item(18, 2).Value = 5
$wsCats.Cells.Item(18, 3).Value = "Universität, Hochschule"
$wsCats.Cells.Item(19, 1).Value = "m_berufsab"
$wsCats.Cells.Item(19, 2).Value = 6
$wsCats.Cells.Item(19, 3).Value = "anderer Abschluss (bis 3/2016 berufl. Praktikum)                          "
$wsCats.Cells.Item(20, 1).Value = "m_berufsab"
$wsCats.Cells.Item(20, 2).Value = 7
$wsCats.Cells.Item(20, 3).Value = "kein beruflicher Abschluss     "
$wsCats.Cells.Item(21, 1).Value = "m_berufsab"
$wsCats.Cells.Item(21, 2).Value = 8
$wsCats.Cells.Item(21, 3).Value = "noch in beruflicher Ausbildung (Auszubildener / Student)"
$wsCats.Cells.Item(22, 1).Value = "v_berufsab"
$wsCats.Cells.Item(22, 2).Value = 1
$wsCats.Cells.Item(22, 3).Value = "Lehre (beruflich-betriebliche Ausbildung)"
$wsCats.Cells.Item(23, 1).Value = "v_berufsab"
$wsCats.Cells.Item(23, 2).Value = 2
$wsCats.Cells.Item(23, 3).Value = "Berufsschule, Handelsschule (berufl.-schulische Ausbildung)"
$wsCats.Cells.Item(24, 1).Value = "v_berufsab"
$wsCats.Cells.Item(24, 2).Value = 3
$wsCats.Cells.Item(24, 3).Value = "Fachschule (z.B. Meister-Technikerschule, Berufs/Fachakademie)"
$wsCats.Cells.Item(25, 1).Value = "v_berufsab"
$wsCats.Cells.Item(25, 2).Value = 4
$wsCats.Cells.Item(25, 3).Value = "Fachhochschule, Ingenieurschule"
$wsCats.Cells.Item(26, 1).Value = "v_berufsab"
$wsCats.Cells.Item(26, 2).Value = 5
$wsCats.Cells.Item(26, 3).Value = "Universität, Hochschule"
$wsCats.Cells.Item(27, 1).Value = "v_berufsab"
$wsCats.Cells.Item(27, 2).Value = 6
$wsCats.Cells.Item(27, 3).Value = "anderer Abschluss (bis 3/2016 berufl. Praktikum)                          "
$wsCats.Cells.Item(28, 1).Value = "v_berufsab"
$wsCats.Cells.Item(28, 2).Value = 7
$wsCats.Cells.Item(28, 3).Value = "kein beruflicher Abschluss     "
$wsCats.Cells.Item(29, 1).Value = "v_berufsab"
$wsCats.Cells.Item(29, 2).Value = 8
$wsCats.Cells.Item(29, 3).Value = "noch in beruflicher Ausbildung (Auszubildener / Student)"

# --- View / selection state ---
$wsVars.Activate()
$wsVars.Range("A2:A57").Select()
$wsCats.Activate()
$wsCats.Range("J38").Select()

Write-Host "done"
